$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 341; everything currently at/after row 341
# (rows 341-352) shifts down to rows 342-353.
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with the new weekly record.
$ws.Cells.Item(341, 1).Value = 10
$ws.Cells.Item(341, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(341, 3).Value = "La Araucanía"
$ws.Cells.Item(341, 4).Value = 45075
$ws.Cells.Item(341, 4).NumberFormat = $ws.Cells.Item(342, 4).NumberFormat
$ws.Cells.Item(341, 5).Value = 9
$ws.Cells.Item(341, 6).Value = 100112043
$ws.Cells.Item(341, 7).Value = "Pepino dulce"
$ws.Cells.Item(341, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 125
$ws.Cells.Item(341, 11).Value = 17000
$ws.Cells.Item(341, 12).Value = 17000
$ws.Cells.Item(341, 13).Value = 17000
$ws.Cells.Item(341, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(341, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(341, 16).Value = 944
$ws.Cells.Item(341, 17).Value = 18
$ws.Cells.Item(341, 18).Value = "Hortaliza"
